$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '30.621.86'
Set-TextValue $ws 'E2' '  +0.69%  '
Set-TextValue $ws 'D3' '1.923.18'
Set-TextValue $ws 'E3' '  -0.26%  '
Set-TextValue $ws 'E4' '  +0.07%  '
Set-TextValue $ws 'D5' '246.96'
Set-TextValue $ws 'E5' '  +2.46%  '
Set-TextValue $ws 'E6' '  +0.04%  '
Set-TextValue $ws 'D7' '0.4744'
Set-TextValue $ws 'E7' '  -0.35%  '
Set-TextValue $ws 'D8' '0.2885'
Set-TextValue $ws 'E8' '  +0.92%  '
Set-TextValue $ws 'D9' '0.06817'
Set-TextValue $ws 'E9' '  +3.15%  '
Set-TextValue $ws 'D10' '105.26'
Set-TextValue $ws 'E10' '  -0.17%  '
Set-TextValue $ws 'D11' '18.38'
Set-TextValue $ws 'E11' '  -4.08%  '
Set-TextValue $ws 'D12' '1.923.37'
Set-TextValue $ws 'E12' '  +0.19%  '
Set-TextValue $ws 'D13' '0.07693'
Set-TextValue $ws 'E13' '  +1.20%  '
Set-TextValue $ws 'D14' '5.344'
Set-TextValue $ws 'E14' '  +4.26%  '
Set-TextValue $ws 'D15' '0.6677'
Set-TextValue $ws 'E15' '  +1.43%  '
Set-TextValue $ws 'D16' '291.64'
Set-TextValue $ws 'E16' '  -3.20%  '
Set-TextValue $ws 'D17' '30.616.74'
Set-TextValue $ws 'E17' '  +0.69%  '
Set-TextValue $ws 'D18' '0.000007614'
Set-TextValue $ws 'E18' '  +1.27%  '
Set-TextValue $ws 'E19' '  +0.01%  '
Set-TextValue $ws 'E20' '  +0.30%  '
Set-TextValue $ws 'D21' '5.546'
Set-TextValue $ws 'E21' '  +5.55%  '
Set-TextValue $ws 'D22' '2.171.28'
Set-TextValue $ws 'E22' '  +0.44%  '
Set-TextValue $ws 'D23' '0.9999'
Set-TextValue $ws 'E23' '  +0.03%  '
Set-TextValue $ws 'D24' '6.444'
Set-TextValue $ws 'E24' '  +1.94%  '
Set-TextValue $ws 'D25' '9.468'
Set-TextValue $ws 'E25' '  +2.53%  '
Set-TextValue $ws 'D26' '167.57'
Set-TextValue $ws 'E26' '  -0.52%  '
Set-TextValue $ws 'D27' '21.11'
Set-TextValue $ws 'E27' '  +6.98%  '
Set-TextValue $ws 'D28' '2.119'
Set-TextValue $ws 'E28' '  +5.39%  '
Set-TextValue $ws 'E29' '  -5.09%  '
Set-TextValue $ws 'D30' '1.400'
Set-TextValue $ws 'E30' '  +3.55%  '
Set-TextValue $ws 'D31' '4.191'
Set-TextValue $ws 'E31' '  +2.16%  '
Set-TextValue $ws 'D32' '4.057'
Set-TextValue $ws 'E32' '  +3.42%  '
Set-TextValue $ws 'D33' '0.05032'
Set-TextValue $ws 'E33' '  +0.39%  '
Set-TextValue $ws 'D34' '0.7375'
Set-TextValue $ws 'E34' '  -0.65%  '
Set-TextValue $ws 'D35' '1.142'
Set-TextValue $ws 'E35' '  -0.58%  '
Set-TextValue $ws 'D36' '0.02066'
Set-TextValue $ws 'E36' '  +5.84%  '
Set-TextValue $ws 'D37' '2.739'
Set-TextValue $ws 'E37' '  +0.35%  '
Set-TextValue $ws 'D38' '2.688'
Set-TextValue $ws 'E38' '  -0.43%  '
Set-TextValue $ws 'D39' '2.050'
Set-TextValue $ws 'E39' '  +0.23%  '
Set-TextValue $ws 'D40' '111.43'
Set-TextValue $ws 'E40' '  +3.66%  '
Set-TextValue $ws 'D41' '0.8713'
Set-TextValue $ws 'E41' '  -0.29%  '
Set-TextValue $ws 'D42' '0.4388'
Set-TextValue $ws 'E42' '  +6.12%  '
Set-TextValue $ws 'D43' '5.915'
Set-TextValue $ws 'E43' '  +2.08%  '
Set-TextValue $ws 'E44' '  +0.04%  '
Set-TextValue $ws 'D45' '67.86'
Set-TextValue $ws 'E45' '  -3.30%  '
Set-TextValue $ws 'D46' '7.294'
Set-TextValue $ws 'E46' '  +1.07%  '
Set-TextValue $ws 'D47' '9.327'
Set-TextValue $ws 'E47' '  +0.68%  '
Set-TextValue $ws 'D48' '48.24'
Set-TextValue $ws 'E48' '  +15.49%  '
Set-TextValue $ws 'D49' '0.1242'
Set-TextValue $ws 'E49' '  +3.25%  '
Set-TextValue $ws 'B50' 'WOONetwork'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue $ws 'D50' '0.2516'
Set-TextValue $ws 'E50' '  +11.77%  '
Set-TextValue $ws 'B51' 'Elrond'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws 'D51' '34.93'
Set-TextValue $ws 'E51' '  +0.35%  '

Write-Host "Applied all cell updates"
